$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -3
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -5
$ws.Range("F13").Value = -2
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 0
$ws.Range("F20").Value = -5
